{"js": "// Three small text edits on the cover/intro pages:\n//   1) \"Version 1.0 Draft\"  -> \"Version 1.2\"\n//   2) \"07/14/2019\"         -> \"09/25/2019\"\n//   3) \"...implemented as a mobile application.\" -> \"...implemented both a web and mobile application.\"\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${searchText}`);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\nawait replaceOnce(\"Version 1.0 Draft\", \"Version 1.2\");\nawait replaceOnce(\"07/14/2019\", \"09/25/2019\");\nawait replaceOnce(\n  \"It will be implemented as a mobile application.\",\n  \"It will be implemented both a web and mobile application.\"\n);\n", "ps1": "# Three small text edits on the cover/intro pages:\n#   1) \"Version 1.0 Draft\"  -> \"Version 1.2\"\n#   2) \"07/14/2019\"         -> \"09/25/2019\"\n#   3) \"...implemented as a mobile application.\" -> \"...implemented both a web and mobile application.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text \"Version 1.0 Draft\" \"Version 1.2\"\nReplace-Text \"07/14/2019\" \"09/25/2019\"\nReplace-Text \"implemented as a mobile application.\" \"implemented both a web and mobile application.\"\n"}
